$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header values (C1 becomes WIN, new columns D-G, ExpPoints moves to H1)
$ws.Range("C1").Value = "WIN"
$ws.Range("D1").Value = "TOP4"
$ws.Range("E1").Value = "TOP5"
$ws.Range("F1").Value = "TOP6"
$ws.Range("G1").Value = "RELEGATION"
$ws.Range("H1").Value = "ExpPoints"

# Apply the same header formatting (bold, border, centered) used by A1 to the
# newly added header cells C1:H1
$ws.Range("A1").Copy()
$ws.Range("C1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Team names in new order (rows 2-21)
$teams = @(
    "Barcelona",
    "Real Madrid",
    "Atlético de Madrid",
    "Villarreal",
    "Real Betis",
    "Athletic Club",
    "Rayo Vallecano",
    "Espanyol",
    "Celta de Vigo",
    "Sevilla",
    "Getafe",
    "Osasuna",
    "Valencia",
    "Real Sociedad",
    "Elche",
    "Alavés",
    "Mallorca",
    "Levante",
    "Real Oviedo",
    "Girona"
)

$expPoints = @(
    86.16881747842982,
    85.26428057540981,
    72.25908403917991,
    65.36590346162168,
    62.71545418065214,
    55.76817160320594,
    50.50918057558626,
    49.31354835086387,
    48.57103881814587,
    47.346414594265,
    47.32073526818876,
    45.45701112562713,
    44.02633958392919,
    43.661992154039,
    41.66310020443806,
    41.07288562634075,
    40.09252545951451,
    35.35084714632725,
    32.6706848685417,
    32.28616982701777
)

for ($i = 0; $i -lt 20; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $teams[$i]
    $ws.Cells.Item($row, 3).Value = ""
    $ws.Cells.Item($row, 4).Value = ""
    $ws.Cells.Item($row, 5).Value = ""
    $ws.Cells.Item($row, 6).Value = ""
    $ws.Cells.Item($row, 7).Value = ""
    $ws.Cells.Item($row, 8).Value = $expPoints[$i]
}
